# "Edit and first run" - update the exported-plot filename to include the
# date range, move the active selection, and widen column F so the longer
# file name is fully visible (re-creating Excel's own best-fit resize).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# exportFile value (F2): MapaSuicidio.png -> MapaSuicidio2009-2018.png
$ws.Range("F2").Value = "./plots/MapaSuicidio2009-2018.png"

# Move the active cell/selection from E7 to C9 (first run after the edit)
$ws.Range("C9").Select()

# Re-fit column F now that it holds the longer path string
$ws.Columns.Item(6).ColumnWidth = 31.83
